# Update column C ("Förändrad") for all data rows (2..367) from serial date
# 45186 (2023-09-17) to 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = 45186
$newValue = 45188

for ($row = 2; $row -le 367; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
